$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("organ") values -------------------------------------------
# Cell order below is deliberate: it reproduces the exact order in which
# Excel first registers each new distinct string into sharedStrings.xml
# (Uterus before Skin), matching the author's original edit order.
$ws.Range("C4").Value  = "Uterus"
$ws.Range("C2").Value  = "Skin"
$ws.Range("C3").Value  = "Skin"
$ws.Range("C5").Value  = "Small Bowel"
$ws.Range("C6").Value  = "Uterus"
$ws.Range("C7").Value  = "Small Bowel"
$ws.Range("C8").Value  = "Uterus"
$ws.Range("C9").Value  = "Small Bowel"
$ws.Range("C10").Value = "Uterus"
$ws.Range("C11").Value = "Liver"
$ws.Range("C12").Value = "Nerve"
$ws.Range("C13").Value = "Uterus"
$ws.Range("C14").Value = "Soft Tissue"
$ws.Range("C15").Value = "Skin"
$ws.Range("C16").Value = "Blood"
$ws.Range("C17").Value = "Bone"
$ws.Range("C18").Value = "Soft Tissue"

# --- Column A width --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.17

# --- View: zoom + selection -------------------------------------------------
$excel.ActiveWindow.Zoom = 140
$ws.Range("A3").Select()
